# Generate Report for handoff
# The previous handoff attempt for "3f36649c-3c15-464d-9aef-631dbafb71ca.md" failed to
# transform, so a brand new report file name is generated and the status/handoff data
# for that row is reset across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldFileName = "3f36649c-3c15-464d-9aef-631dbafb71ca.md"
$newFileName = "57c10ec3-0993-4e4e-8265-e7182a103845.md"
$newStatus   = "Handoff transform failed"
$emptyDate   = "0001-01-01 00:00:00"

$baseUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/2743af4021f47d2beb1350a90bd80b271e30ce7d"
$newFileUrl  = "$baseUrl/e2e/$newFileName"
$configUrl   = "$baseUrl/.localization-config"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the stale hyperlinks (editing hyperlinks loaded from disk in place is
# unreliable) and rebuild them against the refreshed target file name.
$ov.Range("A1:C3").Hyperlinks.Delete()

$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus

$ov.Hyperlinks.Add($ov.Range("A2"), $newFileUrl, "", "", $newFileName)
$ov.Hyperlinks.Add($ov.Range("A3"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A1:I3").Hyperlinks.Delete()

    $ws.Range("B2").Value = $newStatus
    # The handoff transform failed, so there is no more "Latest Handoff File" -
    # clear the cell (and its old hyperlink) entirely.
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = $emptyDate
    $ws.Range("G2").Value = $emptyDate
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = $emptyDate
    $ws.Range("G3").Value = $emptyDate
    $ws.Range("H3").Value = "Ignored"

    $ws.Hyperlinks.Add($ws.Range("A2"), $newFileUrl, "", "", $newFileName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config")
}

$wb.Save()
